# Reorder the rows on the "Parameters" sheet so that the "foi" (Force of
# infection) parameter - a derived/computed parameter with no independent
# Databook Page / Targetable / Default Value settings - moves from row 4
# down to the bottom of the table (row 7), with the remaining parameters
# ("recrate", "infdeath", "susdeath") each shifting up one row to fill the
# gap. This mirrors dragging row 4 and dropping it below row 7 in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Activate()

$firstCol = 1   # A
$lastCol  = 8   # H
$firstRow = 4
$lastRow  = 7

# 1) Remember row 4's ("foi") current values before anything is overwritten.
$oldVals = @()
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $oldVals += , ($ws.Cells.Item($firstRow, $c).Value())
}

# 2) Shift rows 5-7 up into rows 4-6 (values only - the destination cells
#    already carry matching formatting for every column except G, handled
#    below, so a plain value copy reproduces the row move faithfully).
for ($r = $firstRow; $r -le $lastRow - 1; $r++) {
    $srcRow = $r + 1
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($srcRow, $c).Value()
    }
}

# 2a) Column G carries an explicit percentage number format on some rows
#     only; make sure it travels with its value during the shift.
$ws.Cells.Item($lastRow - 2, 7).NumberFormat = $ws.Cells.Item($lastRow - 1, 7).NumberFormat

# 3) Write the original row 4 ("foi") values into row 7, the new last row.
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $ws.Cells.Item($lastRow, $c).Value = $oldVals[$c - $firstCol]
}

# 3a) "foi" never had a value in column G - make sure the row we vacated
#     that value from doesn't leave stale content/formatting behind.
$ws.Cells.Item($lastRow, 7).Clear()

# 4) Match Excel's natural end state after dragging a whole row: the moved
#    row ends up selected at its destination.
$ws.Range("A4:XFD4").Select()
